# Update EC database: consolidate all rows into a single worker
# (ERIK PICO ECHENIQUE) across periods 2208..2508, update summary
# fields, and remove the now-obsolete trailing worker rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Summary header fields -------------------------------------------------
$ws.Range("E11").Value = 7400000   # VALOR MORA
$ws.Range("C13").Value = 1         # Cant. Trabajadores
$ws.Range("F13").Value = 37        # Cant. Periodos

# ---- Preserve the special "closing" bottom-border formatting --------------
# Row 58 (the last data row before the edit) carries the table's bottom
# border style. The new last data row is row 52, so copy that formatting
# there before the content is rewritten and the extra rows are removed.
$ws.Range("B58:J58").Copy()
$ws.Range("B52:J52").PasteSpecial(-4122)  # xlPasteFormats

# ---- Rewrite the data rows (16-52) for ERIK PICO ECHENIQUE -----------------
$periods = @(
  "2208","2209","2210","2211","2212",
  "2301","2302","2303","2304","2305","2306","2307","2308","2309","2310","2311","2312",
  "2401","2402","2403","2404","2405","2406","2407","2408","2409","2410","2411","2412",
  "2501","2502","2503","2504","2505","2506","2507","2508"
)

$row = 16
foreach ($periodo in $periods) {
    $ws.Range("B$row").Value = "CC"
    $ws.Range("C$row").Value = "1143383102"
    $ws.Range("D$row").Value = "ERIK PICO ECHENIQUE"
    $ws.Range("E$row").Value = $periodo
    $ws.Range("F$row").Value = 200000
    $ws.Range("G$row").Value = 5000000
    $row++
}

# ---- Remove the obsolete rows for the other three workers ------------------
# (previously rows 53-58; the footer rows below shift up automatically)
$ws.Rows("53:58").Delete()
